$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-06 20:47:46"
$ws.Range("E3").Value = "2026-02-06 20:47:49"
$ws.Range("H3").NumberFormat = "@"
$ws.Range("H3").Value = "71%"
$ws.Range("O3").Value = "-2.2 °C"
$ws.Range("E4").Value = "2026-02-06 20:47:51"
$ws.Range("J4").Value = "997.6 hPa"
$ws.Range("O4").Value = "13.6 °C"
$ws.Range("E5").Value = "2026-02-06 20:47:54"
$ws.Range("J5").Value = "997.9 hPa"
$ws.Range("O5").Value = "11.1 °C"
$ws.Range("E6").Value = "2026-02-06 20:47:56"
$ws.Range("J6").Value = "999.0 hPa"
$ws.Range("O6").Value = "15.3 °C"
$ws.Range("E7").Value = "2026-02-06 20:47:59"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "59%"
$ws.Range("J7").Value = "998.7 hPa"
$ws.Range("E8").Value = "2026-02-06 20:48:02"
$ws.Range("O8").Value = "9.9 °C"
$ws.Range("E9").Value = "2026-02-06 20:48:04"
$ws.Range("E10").Value = "2026-02-06 20:48:07"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = "89%"
$ws.Range("E11").Value = "2026-02-06 20:48:09"
$ws.Range("J11").Value = "999.2 hPa"
$ws.Range("E12").Value = "2026-02-06 20:48:11"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "62%"
$ws.Range("O12").Value = "13.5 °C"
$ws.Range("E13").Value = "2026-02-06 20:48:14"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "77%"
$ws.Range("O13").Value = "10.2 °C"
$ws.Range("E14").Value = "2026-02-06 20:48:17"
$ws.Range("O14").Value = "-4.3 °C"
$ws.Range("E15").Value = "2026-02-06 20:48:19"
$ws.Range("O15").Value = "10.5 °C"
$ws.Range("E16").Value = "2026-02-06 20:48:21"
$ws.Range("O16").Value = "6.0 °C"
$ws.Range("E17").Value = "2026-02-06 20:48:24"
$ws.Range("E18").Value = "2026-02-06 20:48:26"
$ws.Range("E19").Value = "2026-02-06 20:48:29"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "77%"
$ws.Range("J19").Value = "1000.2 hPa"
$ws.Range("E20").Value = "2026-02-06 20:48:31"
$ws.Range("O20").Value = "-2.1 °C"
$ws.Range("E21").Value = "2026-02-06 20:48:34"
$ws.Range("K21").Value = "10.2 MJ/m2"
$ws.Range("L21").Value = "23.8 km/h - 232º 20:27 TU"
$ws.Range("E22").Value = "2026-02-06 20:48:37"
$ws.Range("H22").NumberFormat = "@"
$ws.Range("H22").Value = "79%"
$ws.Range("E23").Value = "2026-02-06 20:48:39"
$ws.Range("H23").NumberFormat = "@"
$ws.Range("H23").Value = "84%"
$ws.Range("E24").Value = "2026-02-06 20:48:42"
$ws.Range("J24").Value = "997.5 hPa"
$ws.Range("E25").Value = "2026-02-06 20:48:44"
$ws.Range("J25").Value = "998.8 hPa"
$ws.Range("O25").Value = "4.5 °C"
$ws.Range("E26").Value = "2026-02-06 20:48:47"
$ws.Range("O26").Value = "-1.0 °C"
$ws.Range("E27").Value = "2026-02-06 20:48:49"
$ws.Range("J27").Value = "998.0 hPa"
$ws.Range("E28").Value = "2026-02-06 20:48:52"
$ws.Range("J28").Value = "1000.1 hPa"
$ws.Range("E29").Value = "2026-02-06 20:48:54"
$ws.Range("O29").Value = "12.5 °C"
$ws.Range("E30").Value = "2026-02-06 20:48:57"
$ws.Range("E31").Value = "2026-02-06 20:48:59"
$ws.Range("J31").Value = "999.5 hPa"
$ws.Range("E32").Value = "2026-02-06 20:49:01"
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = "52%"
$ws.Range("O32").Value = "15.4 °C"
$ws.Range("E33").Value = "2026-02-06 20:49:04"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "86%"
$ws.Range("O33").Value = "10.3 °C"
$ws.Range("E34").Value = "2026-02-06 20:49:06"
$ws.Range("E35").Value = "2026-02-06 20:49:09"
$ws.Range("E36").Value = "2026-02-06 20:49:11"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "66%"
$ws.Range("O36").Value = "12.7 °C"
